$d = $word.ActiveDocument

# Wipe the single-paragraph body (this also removes the "_GoBack" bookmark
# and the spell-check proofErr markers that bracket "git").
$target = $d.Paragraphs(1).Range
$target.Delete()

# Rebuild the body as three paragraphs:
#   1) "Just for test, start with experiment with git" (single merged run)
#   2) an empty paragraph
#   3) "Added new line ----exper." followed by the relocated "_GoBack" bookmark
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p><w:r><w:t>Just for test, start with experiment with git</w:t></w:r></w:p>
          <w:p/>
          <w:p><w:r><w:t>Added new line ----exper.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$d.Paragraphs(1).Range.InsertXML($xml)
